# Renaming conditionsCombination to conditionsPermutation
#
# The only user-visible text changes in the target revision are:
#   1) "...the conditions combination, which I'll call..."
#        -> "...the conditions permutation, which I'll call..."
#   2) "conditionsCombinationAndInfoGain" -> "conditionsPermutationAndInfoGain"
#
# (Everything else in the underlying OOXML diff -- proofErr spell-check tags,
#  run splitting/merging, namespace/mc:Ignorable bookkeeping, the hidden
#  "_GoBack" bookmark, and styles.xml lsdException bookkeeping -- carries no
#  visible text change, so we only need these two precise replacements.)

$d = $word.ActiveDocument

# 1) "conditions combination" -> "conditions permutation" (lowercase, prose sentence)
$d.Content.Find.Execute(
    "conditions combination",  # FindText
    $false,                    # MatchCase
    $true,                     # MatchWholeWord
    $false,                    # MatchWildcards
    $false,                    # MatchSoundsLike
    $false,                    # MatchAllWordForms
    $true,                     # Forward
    1,                         # Wrap (wdFindContinue)
    $false,                    # Format
    "conditions permutation",  # ReplaceWith
    2                          # Replace (wdReplaceAll)
) | Out-Null

# 2) "conditionsCombinationAndInfoGain" -> "conditionsPermutationAndInfoGain" (identifier)
$d.Content.Find.Execute(
    "conditionsCombinationAndInfoGain",
    $true,                     # MatchCase
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "conditionsPermutationAndInfoGain",
    2
) | Out-Null

Write-Host "done"
